$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rebuild rows 10-23 (the course objectives/program/evaluation block). ---
# The previous layout (rows 10-25) had its rows shifted/remapped, plus two
# trailing rows disappear entirely. Rather than patch each cell in place,
# delete the whole block and insert a clean set of rows so no stale values,
# styles or row heights leak through from the old layout.
$ws.Range("A10:C25").EntireRow.Delete()
$ws.Rows.Item("10:23").Insert()

# Template cells on row 3 already carry the 3 column styles used throughout
# the sheet: col A = bold/top, col B = top/wrap, col C = red/top/wrap.
$styleA = $ws.Range("A3")
$styleB = $ws.Range("B3")
$styleC = $ws.Range("C3")

# Row 10
$ws.Cells.Item(10,1).NumberFormat = "@"
$ws.Cells.Item(10,1).Value = 'Objetivos:'
$ws.Cells.Item(10,2).NumberFormat = "@"
$ws.Cells.Item(10,2).Value = '5983729 - Fernando Vernilli Junior'
$ws.Cells.Item(10,3).NumberFormat = "@"
$ws.Cells.Item(10,3).Value = '5983729 - Fernando Vernilli Junior'
$styleA.Copy()
$ws.Range("A10").PasteSpecial(-4122)
$styleB.Copy()
$ws.Range("B10").PasteSpecial(-4122)
$styleC.Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Rows.Item(10).RowHeight = 60

# Row 11
$ws.Cells.Item(11,1).NumberFormat = "@"
$ws.Cells.Item(11,1).Value = 'Objectives:'
$styleA.Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Rows.Item(11).RowHeight = 60

# Row 12
$ws.Cells.Item(12,1).NumberFormat = "@"
$ws.Cells.Item(12,1).Value = 'Docentes responsáveis:'
$styleA.Copy()
$ws.Range("A12").PasteSpecial(-4122)

# Row 13
$ws.Cells.Item(13,1).NumberFormat = "@"
$ws.Cells.Item(13,1).Value = 'Programa resumido:'
$ws.Cells.Item(13,2).NumberFormat = "@"
$ws.Cells.Item(13,2).Value = '01/01/2022'
$ws.Cells.Item(13,3).NumberFormat = "@"
$ws.Cells.Item(13,3).Value = '01/01/2022'
$styleA.Copy()
$ws.Range("A13").PasteSpecial(-4122)
$styleB.Copy()
$ws.Range("B13").PasteSpecial(-4122)
$styleC.Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Cells.Item(14,1).NumberFormat = "@"
$ws.Cells.Item(14,1).Value = 'Short syllabus:'
$styleA.Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Cells.Item(15,1).NumberFormat = "@"
$ws.Cells.Item(15,1).Value = 'Programa:'
$ws.Cells.Item(15,2).NumberFormat = "@"
$ws.Cells.Item(15,2).Value = '5983729 - Fernando Vernilli Junior'
$ws.Cells.Item(15,3).NumberFormat = "@"
$ws.Cells.Item(15,3).Value = '5983729 - Fernando Vernilli Junior'
$styleA.Copy()
$ws.Range("A15").PasteSpecial(-4122)
$styleB.Copy()
$ws.Range("B15").PasteSpecial(-4122)
$styleC.Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Rows.Item(15).RowHeight = 120

# Row 16
$ws.Cells.Item(16,1).NumberFormat = "@"
$ws.Cells.Item(16,1).Value = 'Syllabus:'
$styleA.Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Cells.Item(17,1).NumberFormat = "@"
$ws.Cells.Item(17,1).Value = 'Avaliação:'
$styleA.Copy()
$ws.Range("A17").PasteSpecial(-4122)

# Row 18
$ws.Cells.Item(18,1).NumberFormat = "@"
$ws.Cells.Item(18,1).Value = 'Método:'
$ws.Cells.Item(18,2).NumberFormat = "@"
$ws.Cells.Item(18,2).Value = '1922320 - Sebastiao Ribeiro'
$ws.Cells.Item(18,3).NumberFormat = "@"
$ws.Cells.Item(18,3).Value = '1922320 - Sebastiao Ribeiro'
$styleA.Copy()
$ws.Range("A18").PasteSpecial(-4122)
$styleB.Copy()
$ws.Range("B18").PasteSpecial(-4122)
$styleC.Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Cells.Item(19,1).NumberFormat = "@"
$ws.Cells.Item(19,1).Value = 'Critério:'
$ws.Cells.Item(19,2).NumberFormat = "@"
$ws.Cells.Item(19,2).Value = 'a) Duas provas escritas (P1 e P2, com peso 1)b) Relatórios sobre os testes experimentais: soma das notas dos relatórios divido pelo número de relatórios (SR), com peso 1.'
$ws.Cells.Item(19,3).NumberFormat = "@"
$ws.Cells.Item(19,3).Value = 'a) Duas provas escritas (P1 e P2, com peso 1)b) Relatórios sobre os testes experimentais: soma das notas dos relatórios divido pelo número de relatórios (SR), com peso 1.'
$styleA.Copy()
$ws.Range("A19").PasteSpecial(-4122)
$styleB.Copy()
$ws.Range("B19").PasteSpecial(-4122)
$styleC.Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Cells.Item(20,1).NumberFormat = "@"
$ws.Cells.Item(20,1).Value = 'Norma de recuperação:'
$ws.Cells.Item(20,2).NumberFormat = "@"
$ws.Cells.Item(20,2).Value = 'Serão aplicadas duas avaliações escritas (P1 e P2, com peso 1). A nota final serão calculada pela equaçãoNF = (P1+P2+MR)/3. NF igual ou superior a 5: aprovação direta. NF entre 3 e 4,9: recuperação. NF inferior a 3: reprovação direta.'
$ws.Cells.Item(20,3).NumberFormat = "@"
$ws.Cells.Item(20,3).Value = 'Serão aplicadas duas avaliações escritas (P1 e P2, com peso 1). A nota final serão calculada pela equaçãoNF = (P1+P2+MR)/3. NF igual ou superior a 5: aprovação direta. NF entre 3 e 4,9: recuperação. NF inferior a 3: reprovação direta.'
$styleA.Copy()
$ws.Range("A20").PasteSpecial(-4122)
$styleB.Copy()
$ws.Range("B20").PasteSpecial(-4122)
$styleC.Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Cells.Item(21,1).NumberFormat = "@"
$ws.Cells.Item(21,1).Value = 'Bibliografia:'
$ws.Cells.Item(21,2).NumberFormat = "@"
$ws.Cells.Item(21,2).Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2. Média final igual ou superior a 5 (cinco): aprovado. NF inferior a 5: reprovado.'
$ws.Cells.Item(21,3).NumberFormat = "@"
$ws.Cells.Item(21,3).Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2. Média final igual ou superior a 5 (cinco): aprovado. NF inferior a 5: reprovado.'
$styleA.Copy()
$ws.Range("A21").PasteSpecial(-4122)
$styleB.Copy()
$ws.Range("B21").PasteSpecial(-4122)
$styleC.Copy()
$ws.Range("C21").PasteSpecial(-4122)
$ws.Rows.Item(21).RowHeight = 120

# Row 22
$ws.Cells.Item(22,1).NumberFormat = "@"
$ws.Cells.Item(22,1).Value = 'Requisitos:'
$styleA.Copy()
$ws.Range("A22").PasteSpecial(-4122)

# Row 23
$ws.Cells.Item(23,2).NumberFormat = "@"
$ws.Cells.Item(23,2).Value = 'LOM3113 -  Tratamentos de Minérios e Hidrometalurgia  (Requisito fraco)
'
$ws.Cells.Item(23,3).NumberFormat = "@"
$ws.Cells.Item(23,3).Value = 'LOM3113 -  Tratamentos de Minérios e Hidrometalurgia  (Requisito fraco)
'
$styleB.Copy()
$ws.Range("B23").PasteSpecial(-4122)
$styleC.Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Rows.Item(23).RowHeight = 30

$excel.CutCopyMode = $false
Write-Host "done"